$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.032835588569617
$ws.Range("D2").Value = 1.049746824815685
$ws.Range("E2").Value = 1.04317353377734
$ws.Range("F2").Value = 1.056414621338521
$ws.Range("I2").Value = 1.041650676869526
$ws.Range("J2").Value = 1.037963623736653
$ws.Range("K2").Value = 1.052503024023049
$ws.Range("L2").Value = 1.045948145618356
$ws.Range("M2").Value = 1.059152396843716
$ws.Range("N2").Value = 1.039437651138451
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.033894147830388
$ws.Range("D3").Value = 1.050374660296304
$ws.Range("E3").Value = 1.044063665628472
$ws.Range("F3").Value = 1.057270597540094
$ws.Range("I3").Value = 1.041801286884613
$ws.Range("J3").Value = 1.03866409079077
$ws.Range("K3").Value = 1.05294364372757
$ws.Range("L3").Value = 1.046649051591924
$ws.Range("M3").Value = 1.059821899589058
$ws.Range("N3").Value = 1.04013911293613
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.034579471716971
$ws.Range("D4").Value = 1.050779786906604
$ws.Range("E4").Value = 1.044639957484732
$ws.Range("F4").Value = 1.057823972856086
$ws.Range("I4").Value = 1.041896541790686
$ws.Range("J4").Value = 1.039117144318417
$ws.Range("K4").Value = 1.053226940795282
$ws.Range("L4").Value = 1.047102267197189
$ws.Range("M4").Value = 1.06025397154791
$ws.Range("N4").Value = 1.040592809851752
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.034867669129139
$ws.Range("D5").Value = 1.050949831068762
$ws.Range("E5").Value = 1.044882305482785
$ws.Range("F5").Value = 1.058056491055373
$ws.Range("I5").Value = 1.041936059390752
$ws.Range("J5").Value = 1.039307560662602
$ws.Range("K5").Value = 1.053345603356552
$ws.Range("L5").Value = 1.047292722230212
$ws.Range("M5").Value = 1.060435340145043
$ws.Range("N5").Value = 1.040783496608986
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.03491606385354
$ws.Range("D6").Value = 1.050978366285047
$ws.Range("E6").Value = 1.044923001141246
$ws.Range("F6").Value = 1.058095524759334
$ws.Range("I6").Value = 1.041942663606329
$ws.Range("J6").Value = 1.039339529611255
$ws.Range("K6").Value = 1.053365501742005
$ws.Range("L6").Value = 1.047324695940464
$ws.Range("M6").Value = 1.060465776599052
$ws.Range("N6").Value = 1.040815510957214
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.034583322282312
$ws.Range("D7").Value = 1.050782060112376
$ws.Range("E7").Value = 1.044643195457503
$ws.Range("F7").Value = 1.057827080252087
$ws.Range("I7").Value = 1.041897071900932
$ws.Range("J7").Value = 1.039119688857729
$ws.Range("K7").Value = 1.053228528083368
$ws.Range("L7").Value = 1.047104812369507
$ws.Range("M7").Value = 1.060256396084045
$ws.Range("N7").Value = 1.040595358004601
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.033193257882113
$ws.Range("D8").Value = 1.049959236630396
$ws.Range("E8").Value = 1.043474291765348
$ws.Range("F8").Value = 1.056704004867245
$ws.Range("I8").Value = 1.041702031175077
$ws.Range("J8").Value = 1.038200390002375
$ws.Range("K8").Value = 1.052652308139101
$ws.Range("L8").Value = 1.046185085311781
$ws.Range("M8").Value = 1.059378893937496
$ws.Range("N8").Value = 1.039674753639428
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.030746592271951
$ws.Range("D9").Value = 1.048500777435971
$ws.Range("E9").Value = 1.041417005528034
$ws.Range("F9").Value = 1.054721246979802
$ws.Range("I9").Value = 1.041341536082882
$ws.Range("J9").Value = 1.036578994807443
$ws.Range("K9").Value = 1.051623108156873
$ws.Range("L9").Value = 1.044562013009258
$ws.Range("M9").Value = 1.057823939183687
$ws.Range("N9").Value = 1.038051055877342
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.029117375850761
$ws.Range("D10").Value = 1.047522831413348
$ws.Range("E10").Value = 1.040047199412959
$ws.Range("F10").Value = 1.053396965355843
$ws.Range("I10").Value = 1.041089957395523
$ws.Range("J10").Value = 1.035497104737656
$ws.Range("K10").Value = 1.050927759212578
$ws.Range("L10").Value = 1.043478399869874
$ws.Range("M10").Value = 1.05678153682409
$ws.Range("N10").Value = 1.036967629399562
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.028412357032458
$ws.Range("D11").Value = 1.047098053087122
$ws.Range("E11").Value = 1.039454476986269
$ws.Range("F11").Value = 1.052822973915772
$ws.Range("I11").Value = 1.04097836368793
$ws.Range("J11").Value = 1.035028412515204
$ws.Range("K11").Value = 1.050624496090768
$ws.Range("L11").Value = 1.043008821285832
$ws.Range("M11").Value = 1.056328812165794
$ws.Range("N11").Value = 1.036498271580394
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.02815054809265
$ws.Range("D12").Value = 1.046940074616499
$ws.Range("E12").Value = 1.039234376321035
$ws.Range("F12").Value = 1.052609683745789
$ws.Range("I12").Value = 1.040936514211476
$ws.Range("J12").Value = 1.034854285692063
$ws.Range("K12").Value = 1.050511525537201
$ws.Range("L12").Value = 1.042834344421118
$ws.Range("M12").Value = 1.056160447213367
$ws.Range("N12").Value = 1.036323897477191
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028206703993406
$ws.Range("D13").Value = 1.046973970418531
$ws.Range("E13").Value = 1.039281585804088
$ws.Range("F13").Value = 1.052655439001108
$ws.Range("I13").Value = 1.040945509091192
$ws.Range("J13").Value = 1.03489163800377
$ws.Range("K13").Value = 1.050535772797109
$ws.Range("L13").Value = 1.042871772755279
$ws.Range("M13").Value = 1.056196571229346
$ws.Range("N13").Value = 1.036361302833466
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.028390714468819
$ws.Range("D14").Value = 1.047084998544621
$ws.Range("E14").Value = 1.039436282109987
$ws.Range("F14").Value = 1.052805344994058
$ws.Range("I14").Value = 1.040974912523974
$ws.Range("J14").Value = 1.035014019822142
$ws.Range("K14").Value = 1.050615164534857
$ws.Range("L14").Value = 1.042994400081165
$ws.Range("M14").Value = 1.056314899198202
$ws.Range("N14").Value = 1.036483858448057
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.028504098197965
$ws.Range("D15").Value = 1.047153380596046
$ws.Range("E15").Value = 1.039531603960855
$ws.Range("F15").Value = 1.052897695904932
$ws.Range("I15").Value = 1.040992976154123
$ws.Range("J15").Value = 1.035089418850553
$ws.Range("K15").Value = 1.050664037359536
$ws.Range("L15").Value = 1.043069947625561
$ws.Range("M15").Value = 1.056387778117346
$ws.Range("N15").Value = 1.036559364551737
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029164174948516
$ws.Range("D16").Value = 1.047550994853063
$ws.Range("E16").Value = 1.040086545214257
$ws.Range("F16").Value = 1.053435047401878
$ws.Range("I16").Value = 1.041097307559982
$ws.Range("J16").Value = 1.035528205554573
$ws.Range("K16").Value = 1.050947840124074
$ws.Range("L16").Value = 1.043509556579252
$ws.Range("M16").Value = 1.05681155416074
$ws.Range("N16").Value = 1.036998774383206
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.029578342119964
$ws.Range("D17").Value = 1.047800055183062
$ws.Range("E17").Value = 1.040434756275619
$ws.Range("F17").Value = 1.053771962418651
$ws.Range("I17").Value = 1.04116204071421
$ws.Range("J17").Value = 1.035803384451013
$ws.Range("K17").Value = 1.051125281369244
$ws.Range("L17").Value = 1.043785213887395
$ws.Range("M17").Value = 1.057077015001495
$ws.Range("N17").Value = 1.037274344065243
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.029819961536754
$ws.Range("D18").Value = 1.04794520027152
$ws.Range("E18").Value = 1.040637901656428
$ws.Range("F18").Value = 1.053968424287321
$ws.Range("I18").Value = 1.041199541910394
$ws.Range("J18").Value = 1.035963869747878
$ws.Range("K18").Value = 1.051228570153091
$ws.Range("L18").Value = 1.043945964682909
$ws.Range("M18").Value = 1.057231722760912
$ws.Range("N18").Value = 1.037435057269636
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.02990235476059
$ws.Range("D19").Value = 1.047994669305206
$ws.Range("E19").Value = 1.04070717569994
$ws.Range("F19").Value = 1.054035403346812
$ws.Range("I19").Value = 1.041212285308037
$ws.Range("J19").Value = 1.036018587360534
$ws.Range("K19").Value = 1.051263753353203
$ws.Range("L19").Value = 1.044000770517184
$ws.Range("M19").Value = 1.057284451845515
$ws.Range("N19").Value = 1.037489852587579
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.029533901484431
$ws.Range("D20").Value = 1.047773346552022
$ws.Range("E20").Value = 1.040397392405465
$ws.Range("F20").Value = 1.053735820301505
$ws.Range("I20").Value = 1.041155121990481
$ws.Range("J20").Value = 1.03577386262175
$ws.Range("K20").Value = 1.051106265283213
$ws.Range("L20").Value = 1.043755642128789
$ws.Range("M20").Value = 1.057048547113434
$ws.Range("N20").Value = 1.037244780311595
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.02833652613657
$ws.Range("D21").Value = 1.047052308950956
$ws.Range("E21").Value = 1.039390726162358
$ws.Range("F21").Value = 1.052761203729528
$ws.Range("I21").Value = 1.040966264942943
$ws.Range("J21").Value = 1.034977982353046
$ws.Range("K21").Value = 1.050591794611966
$ws.Range("L21").Value = 1.042958290888496
$ws.Range("M21").Value = 1.056280060151683
$ws.Range("N21").Value = 1.036447769801621
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027584073103969
$ws.Range("D22").Value = 1.046597825651915
$ws.Range("E22").Value = 1.038758159094456
$ws.Range("F22").Value = 1.052147937111072
$ws.Range("I22").Value = 1.040845217205988
$ws.Range("J22").Value = 1.034477386266888
$ws.Range("K22").Value = 1.050266446005723
$ws.Range("L22").Value = 1.042456649688177
$ws.Range("M22").Value = 1.055795708766594
$ws.Range("N22").Value = 1.035946462811601
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027982926280147
$ws.Range("D23").Value = 1.04683886321505
$ws.Range("E23").Value = 1.039093460089332
$ws.Range("F23").Value = 1.052473087046983
$ws.Range("I23").Value = 1.040909605201702
$ws.Range("J23").Value = 1.034742780036692
$ws.Range("K23").Value = 1.050439097289616
$ws.Range("L23").Value = 1.04272260880379
$ws.Range("M23").Value = 1.056052583435183
$ws.Range("N23").Value = 1.036212233471
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.029553982162563
$ws.Range("D24").Value = 1.047785415427068
$ws.Range("E24").Value = 1.040414275407227
$ws.Range("F24").Value = 1.053752151538381
$ws.Range("I24").Value = 1.041158249056679
$ws.Range("J24").Value = 1.035787202333248
$ws.Range("K24").Value = 1.051114858481933
$ws.Range("L24").Value = 1.043769004443112
$ws.Range("M24").Value = 1.057061410931632
$ws.Range("N24").Value = 1.037258138967013
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.031378780712415
$ws.Range("D25").Value = 1.048878825510496
$ws.Range("E25").Value = 1.04194856461715
$ws.Range("F25").Value = 1.055234273970026
$ws.Range("I25").Value = 1.041436719144424
$ws.Range("J25").Value = 1.036998335743581
$ws.Range("K25").Value = 1.051890811184859
$ws.Range("L25").Value = 1.044981895261512
$ws.Range("M25").Value = 1.058226953799151
$ws.Range("N25").Value = 1.038470992325708
